# Updates the cryptos price list: refreshed Price/Volume(1h) figures for most
# rows, plus a reorder+refresh of the Binance-PegBSC-USD / PolygonEcosystemToken
# rows (37 and 38 swap which coin they describe).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column cells that receive new values so Excel
# does not reinterpret numeric-looking strings (e.g. "2.44", "0.0000210") as numbers.
$priceCells = @("D2","D3","D4","D6","D7","D8","D9","D11","D12","D13","D14","D15","D16","D17","D19","D20","D21","D22","D23","D24","D26","D27","D28","D29","D30","D31","D32","D34","D36","D37","D38","D39","D41","D42","D43","D45","D46","D47","D48","D49","D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '96.461.01'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').Value = '3.685.07'
$ws.Range('E3').Value = '  -0.35%  '
$ws.Range('D4').Value = '2.44'
$ws.Range('E4').Value = '  +29.66%  '
$ws.Range('E5').Value = '  -0.08%  '
$ws.Range('D6').Value = '229.01'
$ws.Range('E6').Value = '  -2.76%  '
$ws.Range('D7').Value = '651.48'
$ws.Range('E7').Value = '  +0.31%  '
$ws.Range('D8').Value = '0.439'
$ws.Range('E8').Value = '  +3.12%  '
$ws.Range('D9').Value = '1.16'
$ws.Range('E9').Value = '  +9.81%  '
$ws.Range('E10').Value = '  -0.07%  '
$ws.Range('D11').Value = '3.682.65'
$ws.Range('E11').Value = '  -0.37%  '
$ws.Range('D12').Value = '47.93'
$ws.Range('E12').Value = '  +8.67%  '
$ws.Range('D13').Value = '0.211'
$ws.Range('E13').Value = '  +2.66%  '
$ws.Range('D14').Value = '0.0000303'
$ws.Range('E14').Value = '  -1.15%  '
$ws.Range('D15').Value = '6.67'
$ws.Range('E15').Value = '  -0.42%  '
$ws.Range('D16').Value = '4.369.26'
$ws.Range('E16').Value = '  -0.40%  '
$ws.Range('D17').Value = '96.106.20'
$ws.Range('E17').Value = '  -0.03%  '
$ws.Range('E18').Value = '  +1.56%  '
$ws.Range('D19').Value = '3.689.24'
$ws.Range('E19').Value = '  +0.25%  '
$ws.Range('D20').Value = '19.69'
$ws.Range('E20').Value = '  +6.13%  '
$ws.Range('D21').Value = '13.06'
$ws.Range('E21').Value = '  +0.55%  '
$ws.Range('D22').Value = '0.543'
$ws.Range('E22').Value = '  +7.94%  '
$ws.Range('D23').Value = '531.33'
$ws.Range('E23').Value = '  +2.48%  '
$ws.Range('D24').Value = '3.32'
$ws.Range('E24').Value = '  -1.69%  '
$ws.Range('E25').Value = '  +42.21%  '
$ws.Range('D26').Value = '122.89'
$ws.Range('E26').Value = '  +21.51%  '
$ws.Range('D27').Value = '0.0000210'
$ws.Range('E27').Value = '  +0.53%  '
$ws.Range('D28').Value = '6.87'
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('D29').Value = '3.883.74'
$ws.Range('E29').Value = '  -0.43%  '
$ws.Range('D30').Value = '13.12'
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('D31').Value = '13.43'
$ws.Range('E31').Value = '  +11.46%  '
$ws.Range('D32').Value = '3.02'
$ws.Range('E32').Value = '  +1.01%  '
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('D34').Value = '0.187'
$ws.Range('E34').Value = '  +0.67%  '
$ws.Range('E35').Value = '  -1.27%  '
$ws.Range('D36').Value = '33.09'
$ws.Range('E36').Value = '  +3.07%  '
$ws.Range('B37').Value = 'PolygonEcosystemToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D37').Value = '0.617'
$ws.Range('E37').Value = '  +5.79%  '
$ws.Range('B38').Value = 'Binance-PegBSC-USD'
$ws.Range('C38').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D38').Value = '0.995'
$ws.Range('E38').Value = '  -0.78%  '
$ws.Range('D39').Value = '616.60'
$ws.Range('E39').Value = '  -4.63%  '
$ws.Range('E40').Value = '  -0.02%  '
$ws.Range('D41').Value = '8.54'
$ws.Range('E41').Value = '  -2.69%  '
$ws.Range('D42').Value = '7.20'
$ws.Range('E42').Value = '  +5.79%  '
$ws.Range('D43').Value = '0.507'
$ws.Range('E43').Value = '  +18.12%  '
$ws.Range('E44').Value = '  +3.28%  '
$ws.Range('D45').Value = '0.0510'
$ws.Range('E45').Value = '  +13.57%  '
$ws.Range('D46').Value = '40.43'
$ws.Range('E46').Value = '  -0.64%  '
$ws.Range('D47').Value = '2.01'
$ws.Range('E47').Value = '  -1.69%  '
$ws.Range('D48').Value = '0.965'
$ws.Range('E48').Value = '  +0.95%  '
$ws.Range('D49').Value = '9.10'
$ws.Range('E49').Value = '  +7.90%  '
$ws.Range('D50').Value = '2.30'
$ws.Range('E50').Value = '  +1.40%  '
$ws.Range('E51').Value = '  -0.27%  '
